# Colocar en la plantilla el apellido y nombre del usuario.
# - Cambia la celda D3 (tercer empleado, columna SABADO) de "OFICINA" a
#   "OFICINA,DEFAULT-LIBRE".
# - Agrega una nueva fila 5 con un identificador/sello (timestamp) en A5,
#   almacenado como texto (prefijo de comilla) para que no se interprete
#   como un número.
# - Ajusta la selección activa de la hoja a A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualiza el valor de D3 a la nueva cadena con el sufijo DEFAULT-LIBRE.
$ws.Range("D3").Value = "OFICINA,DEFAULT-LIBRE"

# Agrega la nueva fila con el sello/identificador como texto (no numérico),
# forzando el prefijo de comilla para conservarlo como texto.
$ws.Range("A5").Value = "'1758796831"

# Actualiza la selección/visualización de la hoja.
$ws.Range("A9").Select()
